$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.832.31'
$ws.Range("E2").Value = '  -1.68%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.892.76'
$ws.Range("E3").Value = '  -1.39%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7765'
$ws.Range("E5").Value = '  -4.31%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '244.22'
$ws.Range("E6").Value = '  -0.12%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.000'
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3125'
$ws.Range("E8").Value = '  -3.57%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '25.38'
$ws.Range("E9").Value = '  -6.63%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07208'
$ws.Range("E10").Value = '  +0.50%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08081'
$ws.Range("E11").Value = '  -0.23%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.7674'
$ws.Range("E12").Value = '  -2.36%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.489'
$ws.Range("E13").Value = '  +1.57%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.896.44'
$ws.Range("E14").Value = '  -1.42%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '92.39'
$ws.Range("E15").Value = '  -2.51%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.166'
$ws.Range("E16").Value = '  +2.04%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '29.842.76'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.96'
$ws.Range("E18").Value = '  -2.42%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '243.63'
$ws.Range("E19").Value = '  -3.51%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007768'
$ws.Range("E20").Value = '  -0.81%  '
$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.000'
$ws.Range("E21").Value = '  -0.12%  '
$ws.Range("B22").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.150.85'
$ws.Range("E22").Value = '  -0.53%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.126'
$ws.Range("E23").Value = '  +1.86%  '
$ws.Range("E24").Value = '  +0.03%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1559'
$ws.Range("E25").Value = '  -4.11%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.395'
$ws.Range("E26").Value = '  -1.36%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '162.56'
$ws.Range("E28").Value = '  -2.10%  '
$ws.Range("E29").Value = '  -4.34%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.431'
$ws.Range("E30").Value = '  +4.12%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.551'
$ws.Range("E31").Value = '  +0.74%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.473'
$ws.Range("E32").Value = '  +2.62%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.100'
$ws.Range("E33").Value = '  -1.01%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05519'
$ws.Range("E34").Value = '  -2.18%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.261'
$ws.Range("E35").Value = '  -3.16%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7489'
$ws.Range("E36").Value = '  +0.53%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.001'
$ws.Range("E37").Value = '  +0.01%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.631'
$ws.Range("E38").Value = '  -3.14%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01918'
$ws.Range("E39").Value = '  -1.93%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.777'
$ws.Range("E40").Value = '  -1.47%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.139.16'
$ws.Range("E41").Value = '  +9.69%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '73.54'
$ws.Range("E42").Value = '  -0.27%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4415'
$ws.Range("E43").Value = '  -1.76%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.883'
$ws.Range("E44").Value = '  -1.68%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8499'
$ws.Range("E45").Value = '  -0.59%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.9998'
$ws.Range("E46").Value = '  -0.10%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '103.72'
$ws.Range("E47").Value = '  +0.55%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.891'
$ws.Range("E48").Value = '  -2.36%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.924'
$ws.Range("E49").Value = '  -1.11%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.041'
$ws.Range("E50").Value = '  +11.70%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.464'
$ws.Range("E51").Value = '  -2.57%  '
